# Applies updated loading-percent results for rows 2-25 (columns B,C,E-J,M-O)
# Column D, K, L remain 0 and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(10.12852665559538, 8.110879165973371, 23.28269436308833, 36.52196271966533, 19.35828909065191, 11.85533369209589, 16.55605497724123, 7.285453413798845, 18.5965221777066, 16.44604009235571, 16.78810679960771)
  3 = @(9.539272179910361, 7.809056561021271, 23.25007201181097, 36.45209059074492, 19.32894353661445, 11.8948484608195, 16.649200368167, 7.296801695172481, 18.39985805283274, 16.47866047993556, 16.83985221119672)
  4 = @(9.156458652415964, 7.616458609050131, 23.23476389172599, 36.41897273321174, 19.31940662598353, 11.92119185267974, 16.70998376579892, 7.30414494233784, 18.2808831577797, 16.50033622415592, 16.87593288402915)
  5 = @(8.995224807120984, 7.536222902750048, 23.22971748874903, 36.40794642498118, 19.31765289771135, 11.93244996353665, 16.73565641587494, 7.307232014639789, 18.23289475766023, 16.509583996592, 16.89171618579537)
  6 = @(8.968137597798076, 7.52279661344385, 23.22895162802262, 36.40626488689711, 19.3174904383367, 11.93435093940455, 16.7399738677442, 7.307750344079555, 18.22495758741701, 16.51114465025293, 16.89440212710339)
  7 = @(9.154305315797938, 7.615383497737158, 23.23469100361573, 36.41881401938204, 19.31937434193008, 11.92134156652139, 16.71032634054042, 7.30418619214688, 18.28023390335077, 16.50045926274891, 16.87614137429989)
  8 = @(9.929737404303394, 8.008366016982208, 23.27046897117336, 36.49584644942684, 19.34641129557342, 11.86852622482935, 16.58742588048534, 7.28928854867786, 18.52837237048185, 16.45694622308059, 16.80505237554761)
  9 = @(11.28212130982014, 8.718231060261594, 23.37786960830764, 36.72402033148292, 19.46658787488217, 11.78148450662344, 16.37492887477651, 7.263040998315907, 19.02694192289891, 16.38465439585552, 16.69997713618702)
  10 = @(12.17175888851046, 9.199262715964151, 23.47912141544105, 36.93781652280997, 19.59542989970362, 11.72763286763046, 16.23620423255962, 7.245549171606218, 19.39765371191558, 16.33945101481651, 16.64388908915617)
  11 = @(12.55375391694052, 9.408722036711927, 23.52994079183672, 37.0448671119212, 19.66270045701077, 11.70533151160394, 16.17687706771032, 7.237977384389456, 19.56664490057135, 16.32059630804969, 16.6229916180407)
  12 = @(12.69513396790445, 9.4866525757077, 23.54985968812476, 37.08678917569388, 19.68940232908012, 11.6972027391014, 16.15495551876274, 7.23516529591295, 19.63063726797895, 16.3137015879443, 16.61574466800701)
  13 = @(12.66483090602735, 9.46993111901763, 23.54553995336791, 37.07769934700514, 19.68359729428377, 11.69893934210593, 16.1596524910436, 7.235768478800725, 19.61685622149109, 16.31517559474603, 16.61727574985672)
  14 = @(12.5654509702307, 9.415161441708051, 23.53156606798358, 37.04828846684964, 19.66487273484424, 11.704656411422, 16.17506265353592, 7.237744927523333, 19.57190990265094, 16.32002416520046, 16.62238203664306)
  15 = @(12.50415145466233, 9.381431678450028, 23.52309422332762, 37.03045298300941, 19.65356275800689, 11.70819948283631, 16.18457274248584, 7.238962739179346, 19.54437734482473, 16.3230259642799, 16.6255966486063)
  16 = @(12.14633692828313, 9.185382019075226, 23.47589521204567, 36.93101559063165, 19.59120645419202, 11.72913454474379, 16.24015752762064, 7.246051740104484, 19.38661240450639, 16.34071754830625, 16.64534792000556)
  17 = @(11.92100670038497, 9.062682706472199, 23.44815252346841, 36.87250602954669, 19.55515945829031, 11.74254032796137, 16.27522559382247, 7.250499143072174, 19.28988206133133, 16.35200797519946, 16.65864914876766)
  18 = @(11.78926574281742, 8.991229915151687, 23.43264427775015, 36.8397762859154, 19.53524249076671, 11.75045765204143, 16.29575149270477, 7.25309345440963, 19.23428155859253, 16.35866277752597, 16.66673412868268)
  19 = @(11.74429380966688, 8.966887524217368, 23.42747078367286, 36.82885384613584, 19.52863962911253, 11.75317380322613, 16.30276227079154, 7.253978082214861, 19.21546399145867, 16.36094362143177, 16.66954610284662)
  20 = @(11.9452147283434, 9.075835583917302, 23.45105941756684, 36.878639069269, 19.55891235340906, 11.74109186612991, 16.27145572105437, 7.250021955954774, 19.30017579789942, 16.35078944678946, 16.65718822571259)
  21 = @(12.59473014400592, 9.4312865654441, 23.53565230716725, 37.05688978749339, 19.67033941691503, 11.70296858240759, 16.17052153402595, 7.237162900774494, 19.58511215871112, 16.31859337353, 16.62086408998126)
  22 = @(13.00014639819475, 9.655494021744353, 23.59486623209357, 37.18144314451914, 19.75031099969173, 11.67989647343763, 16.10772862120386, 7.22908030960917, 19.77130578884243, 16.29898004431226, 16.60100983018034)
  23 = @(12.78551556483945, 9.53658340289452, 23.56290669486823, 37.11423804072591, 19.70698109989951, 11.69204164129074, 16.1409517008592, 7.233364793857302, 19.67195021600246, 16.30931750025126, 16.61125012414097)
  24 = @(11.9342771171027, 9.069892003483575, 23.44974383508607, 36.87586349127198, 19.55721315471964, 11.7417460611158, 16.27315894426047, 7.250237575609315, 19.29552196117779, 16.35133983325694, 16.65784734509958)
  25 = @(10.93445088847036, 8.533087416318587, 23.34486146271987, 36.65411592715976, 19.42691571504867, 11.80325981446916, 16.42936127384442, 7.263040998315907, 19.02694192289891, 16.38465439585552, 16.69997713618702)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]

  $bj = New-Object "object[,]" 1,9
  $bj[0,0] = $vals[0]   # B
  $bj[0,1] = $vals[1]   # C
  $bj[0,2] = 0          # D (unchanged)
  $bj[0,3] = $vals[2]   # E
  $bj[0,4] = $vals[3]   # F
  $bj[0,5] = $vals[4]   # G
  $bj[0,6] = $vals[5]   # H
  $bj[0,7] = $vals[6]   # I
  $bj[0,8] = $vals[7]   # J
  $ws.Range("B$($row):J$($row)").Value = $bj

  $mo = New-Object "object[,]" 1,3
  $mo[0,0] = $vals[8]   # M
  $mo[0,1] = $vals[9]   # N
  $mo[0,2] = $vals[10]  # O
  $ws.Range("M$($row):O$($row)").Value = $mo
}
